$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (I1/J1 swap meaning; L1/M1 are brand new columns)
$ws.Range("I1").Value = "py"
$ws.Range("J1").Value = "nAMLOrMDS"
$ws.Range("L1").Value = "nNonBreastSolid"
$ws.Range("M1").Value = "incidenceNonBreastSolid"

# Per-row data: col I/J swap to (py, nAMLOrMDS) order, and new L/M (nNonBreastSolid, incidenceNonBreastSolid) values
$rows = @(
    @{Row=2; Py=801.1083333333332; NAml=9.0; NSolid=3.0; IncSolid=37.4},
    @{Row=3; Py=874.5166666666667; NAml=0.0; NSolid=1.0; IncSolid=11.4},
    @{Row=4; Py=1551.25; NAml=0.0; NSolid=4.0; IncSolid=25.8},
    @{Row=5; Py=1624.25; NAml=0.0; NSolid=3.0; IncSolid=18.5},
    @{Row=6; Py=1551.25; NAml=3.0; NSolid=6.0; IncSolid=38.7},
    @{Row=7; Py=2649.1666666666665; NAml=1.0; NSolid=9.0; IncSolid=34.0},
    @{Row=8; Py=2530.0; NAml=1.0; NSolid=12.0; IncSolid=47.4},
    @{Row=9; Py=1452.0; NAml=2.0; NSolid=4.0; IncSolid=27.5},
    @{Row=10; Py=1479.0; NAml=3.0; NSolid=10.0; IncSolid=67.6},
    @{Row=11; Py=1503.0; NAml=4.0; NSolid=2.0; IncSolid=13.3},
    @{Row=12; Py=1485.0; NAml=2.0; NSolid=10.0; IncSolid=67.3},
    @{Row=13; Py=1419.2999999999997; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=14; Py=1435.8999999999999; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=15; Py=870.0; NAml=$null; NSolid=4.0; IncSolid=46.0},
    @{Row=16; Py=961.8333333333333; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=17; Py=870.0; NAml=$null; NSolid=2.0; IncSolid=23.0},
    @{Row=18; Py=966.6666666666666; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=19; Py=3815.0; NAml=0.0; NSolid=$null; IncSolid=$null},
    @{Row=20; Py=3808.0; NAml=0.0; NSolid=$null; IncSolid=$null},
    @{Row=21; Py=3514.0; NAml=2.0; NSolid=$null; IncSolid=$null},
    @{Row=22; Py=3500.0; NAml=0.0; NSolid=$null; IncSolid=$null},
    @{Row=23; Py=8932.083333333334; NAml=2.0; NSolid=$null; IncSolid=$null},
    @{Row=24; Py=8932.083333333334; NAml=4.0; NSolid=$null; IncSolid=$null},
    @{Row=25; Py=984.0; NAml=0.0; NSolid=3.0; IncSolid=30.5},
    @{Row=26; Py=1044.0; NAml=0.0; NSolid=5.0; IncSolid=47.9},
    @{Row=27; Py=1602.5666666666666; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=28; Py=1589.4666666666665; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=29; Py=1598.1999999999998; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=30; Py=4429.424999999999; NAml=2.0; NSolid=$null; IncSolid=$null},
    @{Row=31; Py=4429.424999999999; NAml=1.0; NSolid=$null; IncSolid=$null},
    @{Row=32; Py=4458.299999999999; NAml=3.0; NSolid=$null; IncSolid=$null},
    @{Row=33; Py=3905.4000000000005; NAml=4.0; NSolid=$null; IncSolid=$null},
    @{Row=34; Py=3896.2000000000003; NAml=10.0; NSolid=$null; IncSolid=$null},
    @{Row=35; Py=3905.4000000000005; NAml=8.0; NSolid=$null; IncSolid=$null},
    @{Row=36; Py=2490.0; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=37; Py=2480.0; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=38; Py=2470.0; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=39; Py=2470.0; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=40; Py=1534.5; NAml=0.0; NSolid=3.0; IncSolid=19.6},
    @{Row=41; Py=1498.3333333333335; NAml=1.0; NSolid=4.0; IncSolid=26.7},
    @{Row=42; Py=2505.208333333333; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=43; Py=4411.458333333333; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=44; Py=5000.0; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=45; Py=4994.791666666666; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=46; Py=2314.333333333333; NAml=0.0; NSolid=5.0; IncSolid=21.6},
    @{Row=47; Py=2281.583333333333; NAml=2.0; NSolid=4.0; IncSolid=17.5},
    @{Row=48; Py=2183.333333333333; NAml=1.0; NSolid=2.0; IncSolid=9.16},
    @{Row=49; Py=9546.625; NAml=7.0; NSolid=$null; IncSolid=$null},
    @{Row=50; Py=9546.625; NAml=7.0; NSolid=$null; IncSolid=$null},
    @{Row=51; Py=8918.25; NAml=8.0; NSolid=$null; IncSolid=$null},
    @{Row=52; Py=9027.5; NAml=8.0; NSolid=$null; IncSolid=$null},
    @{Row=53; Py=3692.416666666667; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=54; Py=3662.916666666667; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=55; Py=1566.5; NAml=1.0; NSolid=$null; IncSolid=$null},
    @{Row=56; Py=1566.5; NAml=0.0; NSolid=$null; IncSolid=$null},
    @{Row=57; Py=11448.0; NAml=13.0; NSolid=50.0; IncSolid=43.7},
    @{Row=58; Py=10972.800000000001; NAml=15.0; NSolid=51.0; IncSolid=46.5},
    @{Row=59; Py=7698.333333333334; NAml=6.0; NSolid=$null; IncSolid=$null},
    @{Row=60; Py=7708.666666666667; NAml=2.0; NSolid=$null; IncSolid=$null},
    @{Row=61; Py=8282.083333333334; NAml=2.0; NSolid=$null; IncSolid=$null},
    @{Row=62; Py=8292.916666666668; NAml=6.0; NSolid=$null; IncSolid=$null},
    @{Row=63; Py=3458.5833333333335; NAml=0.0; NSolid=9.0; IncSolid=26.0},
    @{Row=64; Py=3343.0833333333335; NAml=0.0; NSolid=16.0; IncSolid=47.9},
    @{Row=65; Py=5016.525; NAml=1.0; NSolid=9.0; IncSolid=17.9},
    @{Row=66; Py=5137.849999999999; NAml=2.0; NSolid=10.0; IncSolid=19.5},
    @{Row=67; Py=1780.8; NAml=$null; NSolid=3.0; IncSolid=16.8},
    @{Row=68; Py=2178.3; NAml=$null; NSolid=3.0; IncSolid=13.8},
    @{Row=69; Py=1580.8333333333333; NAml=2.0; NSolid=$null; IncSolid=$null},
    @{Row=70; Py=1545.8333333333333; NAml=3.0; NSolid=$null; IncSolid=$null},
    @{Row=71; Py=4808.5; NAml=2.0; NSolid=9.0; IncSolid=18.7},
    @{Row=72; Py=3908.7500000000005; NAml=0.0; NSolid=4.0; IncSolid=10.2},
    @{Row=73; Py=870.25; NAml=0.0; NSolid=2.0; IncSolid=23.0},
    @{Row=74; Py=$null; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=75; Py=$null; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=76; Py=$null; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=77; Py=1134.3333333333335; NAml=0.0; NSolid=1.0; IncSolid=8.82},
    @{Row=78; Py=1168.9166666666667; NAml=1.0; NSolid=3.0; IncSolid=25.7},
    @{Row=79; Py=4980.0; NAml=3.0; NSolid=$null; IncSolid=$null},
    @{Row=80; Py=5015.0; NAml=1.0; NSolid=$null; IncSolid=$null},
    @{Row=81; Py=2092.7999999999997; NAml=2.0; NSolid=$null; IncSolid=$null},
    @{Row=82; Py=2073.6; NAml=1.0; NSolid=$null; IncSolid=$null},
    @{Row=83; Py=6966.200000000001; NAml=1.0; NSolid=$null; IncSolid=$null},
    @{Row=84; Py=4812.900000000001; NAml=5.0; NSolid=$null; IncSolid=$null},
    @{Row=85; Py=7021.1; NAml=0.0; NSolid=$null; IncSolid=$null},
    @{Row=86; Py=4812.900000000001; NAml=0.0; NSolid=$null; IncSolid=$null},
    @{Row=87; Py=8693.333333333332; NAml=5.0; NSolid=$null; IncSolid=$null},
    @{Row=88; Py=8714.666666666666; NAml=8.0; NSolid=$null; IncSolid=$null},
    @{Row=89; Py=8693.333333333332; NAml=11.0; NSolid=$null; IncSolid=$null},
    @{Row=90; Py=6281.6; NAml=0.0; NSolid=17.0; IncSolid=27.1},
    @{Row=91; Py=6344.0; NAml=0.0; NSolid=14.0; IncSolid=22.1},
    @{Row=92; Py=1920.0; NAml=0.0; NSolid=4.0; IncSolid=20.8},
    @{Row=93; Py=1962.6666666666665; NAml=0.0; NSolid=4.0; IncSolid=20.4},
    @{Row=94; Py=1454.3999999999999; NAml=2.0; NSolid=4.0; IncSolid=27.5},
    @{Row=95; Py=1444.8; NAml=0.0; NSolid=0.0; IncSolid=0.0},
    @{Row=96; Py=6759.0; NAml=$null; NSolid=$null; IncSolid=$null},
    @{Row=97; Py=6678.0; NAml=$null; NSolid=$null; IncSolid=$null}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 9).Value  = $r.Py
    $ws.Cells.Item($r.Row, 10).Value = $r.NAml
    $ws.Cells.Item($r.Row, 12).Value = $r.NSolid
    $ws.Cells.Item($r.Row, 13).Value = $r.IncSolid
}
